$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{
    2  = 45757
    3  = 45759
    4  = 45759
    5  = 45758
    6  = 45759
    7  = 45760
    8  = 45760
    9  = 45757
    10 = 45759
    11 = 45760
    12 = 45759
    13 = 45759
    14 = 45760
    15 = 45760
    16 = 45760
    17 = 45758
    18 = 45759
    19 = 45759
}

foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 1).Value = $dates[$row]
}
